$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.142.09'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.276.67'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.46%  '
$ws.Range("E4").Value = '  -0.49%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '111.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '264.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.652'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.29%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.609'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '46.65'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0935'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.27'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.95%  '
$ws.Range("E13").Value = '  +1.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.26'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.616.14'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.859'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.55%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.283.85'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.189.89'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("E19").Value = '  -0.95%  '
$ws.Range("E20").Value = '  -1.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("E22").Value = '  -1.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.88'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.15%  '
$ws.Range("E26").Value = '  +1.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '41.01'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.47%  '
$ws.Range("E29").Value = '  -1.14%  '
$ws.Range("E30").Value = '  -1.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.05'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.45'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.45%  '
$ws.Range("E33").Value = '  -3.51%  '
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.132'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.03%  '
$ws.Range("E36").Value = '  +3.33%  '
$ws.Range("E37").Value = '  -0.82%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.87'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.87%  '
$ws.Range("E39").Value = '  -3.05%  '
$ws.Range("E40").Value = '  +7.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.22'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.46'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.15%  '
$ws.Range("E43").Value = '  -2.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.08'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("E47").Value = '  +3.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.54'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.18%  '
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '99.93'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.51%  '
$ws.Range("E51").Value = '  +10.57%  '
